# Bug fix in XLSX module: add a new "aux.pmid[]" column (N) to the first sheet,
# with sample values in rows 11 and 12.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, styled like the rest of row 1's header cells.
$ws.Range("N1").Value = "aux.pmid[]"

# New data cells for the added column.
# Order matters for shared-string table layout: "1111;" must be interned
# before "123;321" so the two new rows land on the right shared indices.
$ws.Range("N12").Value = "1111;"
$ws.Range("N11").Value = "123;321"

# Move the active selection, matching the saved view state.
$ws.Range("O18").Select() | Out-Null
